$d = $word.ActiveDocument

# --- Paragraph 1: "hvjiyfygiygiyg" -> three runs "H" | "vjiyfygiygiyg" | " " ---

# 1. Capitalize the leading "h" -> "H" while the paragraph is still a single run,
#    so the in-place text replace merges back cleanly (no stray formatting).
$d.Range(0, 1).Text = "H"

# 2. Split "Hvjiyfygiygiyg" into two runs "H" and "vjiyfygiygiyg" by temporarily
#    inserting a paragraph break right after "H" and then deleting that paragraph
#    mark again. Word re-joins the two paragraphs into one paragraph containing
#    two separate (cleanly formatted) runs instead of recombining the text into
#    a single run.
$d.Range(1, 1).InsertParagraphAfter()
$d.Range(1, 2).Delete()

# 3. Append a trailing " " as a third, separate run using the same split/rejoin
#    trick, this time at the end of the paragraph (just before its paragraph
#    mark).
$p1 = $d.Paragraphs(1)
$endPos = $p1.Range.End - 1
$d.Range($endPos, $endPos).InsertParagraphAfter()
$p2 = $d.Paragraphs(2)
$d.Range($p2.Range.Start, $p2.Range.Start).InsertAfter(" ")
$d.Range($endPos, $endPos + 1).Delete()

# --- New paragraph 2: "Changes broski" ---
$p1 = $d.Paragraphs(1)
$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs(2)
$p2.Range.Text = "Changes broski"
